$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the parameter values: A2 (n_repetition) and C2 (nboot) both become 500
$ws.Range("A2").Value = 500
$ws.Range("C2").Value = 500

# Update the active selection to C2, matching the saved view state
$ws.Range("C2").Select()
